# Regenerate the handback status report:
#  - Drop the stale "b47a0c2b-7634-43d8-91d5-7a00f8ae86fc" entry (row 3) from
#    every sheet, including its hyperlinks.
#  - Refresh the "Correspond Handoff/Handback Datetime" timestamps for the
#    remaining "54b515d0-d368-4f1f-99f3-4f8cf3c71886" entry.

$wb = $excel.ActiveWorkbook

function Remove-RowHyperlinks($ws, $row) {
    $found = $true
    while ($found) {
        $found = $false
        foreach ($h in $ws.Hyperlinks) {
            if ($h.Range.Row -eq $row) {
                $h.Delete()
                $found = $true
                break
            }
        }
    }
}

# Overview sheet (A:C) — no timestamp columns here, just drop row 3.
$wsOverview = $wb.Worksheets.Item("Overview")
Remove-RowHyperlinks $wsOverview 3
$wsOverview.Rows.Item(3).Delete()

# zh-cn sheet — refresh the handoff/handback datetimes, then drop row 3.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-13 00:43:07"
$wsZh.Range("H2").Value = "2016-03-13 00:43:29"
Remove-RowHyperlinks $wsZh 3
$wsZh.Rows.Item(3).Delete()

# de-de sheet — refresh the handoff/handback datetimes, then drop row 3.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-13 00:43:11"
$wsDe.Range("H2").Value = "2016-03-13 00:43:35"
Remove-RowHyperlinks $wsDe 3
$wsDe.Rows.Item(3).Delete()
